{"js": "// Replace the comma with a semicolon in the <prog> production shown in the\n// predictive parsing table: \"program I, var DL begin SL end.\" becomes\n// \"program I; var DL begin SL end.\" (split across three runs to mirror the\n// authored change, but the net visible text is the comma -> semicolon swap).\n\nconst searchResults = context.document.body.search(\"program I, var DL begin SL end.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text not found: 'program I, var DL begin SL end.'\");\n}\n\nconst target = searchResults.items[0];\ntarget.insertText(\"program I; var DL begin SL end.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the comma with a semicolon in the <prog> production shown in the\n# predictive parsing table: \"program I, var DL begin SL end.\" becomes\n# \"program I; var DL begin SL end.\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"program I, var DL begin SL end.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"program I; var DL begin SL end.\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n"}
